$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.024660117703589
$ws.Range("D2").Value = 1.028811018128613
$ws.Range("E2").Value = 1.048127862484221
$ws.Range("F2").Value = 1.052409871202034
$ws.Range("I2").Value = 1.029006745280023
$ws.Range("J2").Value = 1.02983348931217
$ws.Range("K2").Value = 1.031626750732128
$ws.Range("L2").Value = 1.050888573368396
$ws.Range("M2").Value = 1.055158681836927
$ws.Range("N2").Value = 1.013920277353502
$ws.Range("B3").Value = 1.019999999999999
$ws.Range("C3").Value = 1.025874052949465
$ws.Range("D3").Value = 1.029702286629079
$ws.Range("E3").Value = 1.049535383494687
$ws.Range("F3").Value = 1.053925447985052
$ws.Range("I3").Value = 1.02920501264117
$ws.Range("J3").Value = 1.030685059668205
$ws.Range("K3").Value = 1.032325808617118
$ws.Range("L3").Value = 1.052106535967849
$ws.Range("M3").Value = 1.056485295928675
$ws.Range("N3").Value = 1.014207566787423
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.026658668492696
$ws.Range("D4").Value = 1.030278007268106
$ws.Range("E4").Value = 1.050446143585993
$ws.Range("F4").Value = 1.05490616427228
$ws.Range("I4").Value = 1.029331538116388
$ws.Range("J4").Value = 1.031234734344365
$ws.Range("K4").Value = 1.032776528947206
$ws.Range("L4").Value = 1.052894116157382
$ws.Range("M4").Value = 1.057343239023775
$ws.Range("N4").Value = 1.014392865815611
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.026988312889572
$ws.Range("D5").Value = 1.030519804606107
$ws.Range("E5").Value = 1.050829031631406
$ws.Range("F5").Value = 1.055318471093627
$ws.Range("I5").Value = 1.029384306854236
$ws.Range("J5").Value = 1.031465496750231
$ws.Range("K5").Value = 1.032965626039096
$ws.Range("L5").Value = 1.053225093705747
$ws.Range("M5").Value = 1.057703811293473
$ws.Range("N5").Value = 1.014470623303644
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.027043649528601
$ws.Range("D6").Value = 1.030560389664215
$ws.Range("E6").Value = 1.050893320614604
$ws.Range("F6").Value = 1.055387700168191
$ws.Range("I6").Value = 1.029393142189077
$ws.Range("J6").Value = 1.031504224016656
$ws.Range("K6").Value = 1.03299735369211
$ws.Range("L6").Value = 1.05328065930515
$ws.Range("M6").Value = 1.057764346813889
$ws.Range("N6").Value = 1.014483670804373
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.026663074030615
$ws.Range("D7").Value = 1.03028123910019
$ws.Range("E7").Value = 1.050451259731461
$ws.Range("F7").Value = 1.054911673476311
$ws.Range("I7").Value = 1.029332244874787
$ws.Range("J7").Value = 1.031237819061155
$ws.Range("K7").Value = 1.032779057185834
$ws.Range("L7").Value = 1.052898539166628
$ws.Range("M7").Value = 1.057348057423217
$ws.Range("N7").Value = 1.014393905373115
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.025070557249288
$ws.Range("D8").Value = 1.029112432204835
$ws.Range("E8").Value = 1.048603543555407
$ws.Range("F8").Value = 1.052922062807437
$ws.Range("I8").Value = 1.029074116723073
$ws.Range("J8").Value = 1.030121561852994
$ws.Range("K8").Value = 1.031863336214158
$ws.Range("L8").Value = 1.051300300278118
$ws.Range("M8").Value = 1.055607116137233
$ws.Range("N8").Value = 1.014017492003685
$ws.Range("B9").Value = 1.019999999999999
$ws.Range("C9").Value = 1.02225744618424
$ws.Range("D9").Value = 1.02704520561646
$ws.Range("E9").Value = 1.045347427534836
$ws.Range("F9").Value = 1.049416176891361
$ws.Range("I9").Value = 1.028605712475588
$ws.Range("J9").Value = 1.028144159500711
$ws.Range("K9").Value = 1.030237274492223
$ws.Range("L9").Value = 1.048479811306144
$ws.Range("M9").Value = 1.052535601645743
$ws.Range("N9").Value = 1.013349609767987
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.020377206064007
$ws.Range("D10").Value = 1.025661831472289
$ws.Range("E10").Value = 1.043176261001814
$ws.Range("F10").Value = 1.047078653458054
$ws.Range("I10").Value = 1.028284306255378
$ws.Range("J10").Value = 1.026818763226645
$ws.Range("K10").Value = 1.029144778362324
$ws.Range("L10").Value = 1.046596414068437
$ws.Range("M10").Value = 1.050485130747561
$ws.Range("N10").Value = 1.012901229327748
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.019561850176292
$ws.Range("D11").Value = 1.025061556074856
$ws.Range("E11").Value = 1.042235957805539
$ws.Range("F11").Value = 1.04606634789478
$ws.Range("I11").Value = 1.028142960273121
$ws.Range("J11").Value = 1.026243132544616
$ws.Range("K11").Value = 1.028669689641255
$ws.Range("L11").Value = 1.045780096537351
$ws.Range("M11").Value = 1.049596530630083
$ws.Range("N11").Value = 1.012706325729597
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.019258806513415
$ws.Range("D12").Value = 1.024838395111485
$ws.Range("E12").Value = 1.041886655659092
$ws.Range("F12").Value = 1.045690305119644
$ws.Range("I12").Value = 1.028090130695426
$ws.Range("J12").Value = 1.026029055884109
$ws.Range("K12").Value = 1.028492913582983
$ws.Range("L12").Value = 1.045476755556809
$ws.Range("M12").Value = 1.049266349599226
$ws.Range("N12").Value = 1.012633816110117
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.019323818806443
$ws.Range("D13").Value = 1.024886272612064
$ws.Range("E13").Value = 1.041961583678134
$ws.Range("F13").Value = 1.045770968908818
$ws.Range("I13").Value = 1.028101477644165
$ws.Range("J13").Value = 1.026074987942913
$ws.Range("K13").Value = 1.028530846568017
$ws.Range("L13").Value = 1.045541828932322
$ws.Range("M13").Value = 1.049337179857892
$ws.Range("N13").Value = 1.012649374829425
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.019536804272895
$ws.Range("D14").Value = 1.02504311344664
$ws.Range("E14").Value = 1.042207085059543
$ws.Range("F14").Value = 1.046035264675809
$ws.Range("I14").Value = 1.028138600045956
$ws.Range("J14").Value = 1.026225442256616
$ws.Range("K14").Value = 1.028655083549755
$ws.Range("L14").Value = 1.045755024842006
$ws.Range("M14").Value = 1.049569240147603
$ws.Range("N14").Value = 1.012700334386728
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.019668007123071
$ws.Range("D15").Value = 1.025139722769497
$ws.Range("E15").Value = 1.042358342199348
$ws.Range("F15").Value = 1.046198102207624
$ws.Range("I15").Value = 1.028161428982282
$ws.Range("J15").Value = 1.026318107356971
$ws.Range("K15").Value = 1.028731589324053
$ws.Range("L15").Value = 1.045886365257058
$ws.Range("M15").Value = 1.049712204702687
$ws.Range("N15").Value = 1.012731717156107
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.020431292938378
$ws.Range("D16").Value = 1.025701642957313
$ws.Range("E16").Value = 1.043238661692622
$ws.Range("F16").Value = 1.047145833329992
$ws.Range("I16").Value = 1.028293641045692
$ws.Range("J16").Value = 1.026856929368209
$ws.Range("K16").Value = 1.029176265492925
$ws.Range("L16").Value = 1.046650573217602
$ws.Range("M16").Value = 1.050544088321745
$ws.Range("N16").Value = 1.012914148529451
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.020909758026191
$ws.Range("D17").Value = 1.026053780589336
$ws.Range("E17").Value = 1.043790812409751
$ws.Range("F17").Value = 1.047740277641588
$ws.Range("I17").Value = 1.028375991442199
$ws.Range("J17").Value = 1.027194454438621
$ws.Range("K17").Value = 1.029454654120227
$ws.Range("L17").Value = 1.047129724555975
$ws.Range("M17").Value = 1.051065707039028
$ws.Range("N17").Value = 1.013028381054772
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.021188723152185
$ws.Range("D18").Value = 1.026259054560626
$ws.Range("E18").Value = 1.044112856194165
$ws.Range("F18").Value = 1.048086993611058
$ws.Range("I18").Value = 1.028423815225197
$ws.Range("J18").Value = 1.027391160806771
$ws.Range("K18").Value = 1.029616837692516
$ws.Range("L18").Value = 1.047409129145539
$ws.Range("M18").Value = 1.051369888078337
$ws.Range("N18").Value = 1.013094938459209
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.021283823585546
$ws.Range("D19").Value = 1.026329027050551
$ws.Range("E19").Value = 1.04422266219042
$ws.Range("F19").Value = 1.048205212714944
$ws.Range("I19").Value = 1.02844008633044
$ws.Range("J19").Value = 1.02745820446589
$ws.Range("K19").Value = 1.029672104943544
$ws.Range("L19").Value = 1.047504386102356
$ws.Range("M19").Value = 1.051473594181983
$ws.Range("N19").Value = 1.013117620541057
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.020858435223541
$ws.Range("D20").Value = 1.026016012182473
$ws.Range("E20").Value = 1.043731573600124
$ws.Range("F20").Value = 1.047676500817689
$ws.Range("I20").Value = 1.028367177730507
$ws.Range("J20").Value = 1.027158258413208
$ws.Range("K20").Value = 1.02942480592244
$ws.Range("L20").Value = 1.047078324079515
$ws.Range("M20").Value = 1.051009749602003
$ws.Range("N20").Value = 1.013016132490355
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.019474090477882
$ws.Range("D21").Value = 1.024996933043915
$ws.Range("E21").Value = 1.042134791940825
$ws.Range("F21").Value = 1.045957436965036
$ws.Range("I21").Value = 1.028127677466379
$ws.Range("J21").Value = 1.026181144458856
$ws.Range("K21").Value = 1.028618507343658
$ws.Range("L21").Value = 1.045692247431998
$ws.Range("M21").Value = 1.049500907404071
$ws.Range("N21").Value = 1.012685331213449
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.018602630386506
$ws.Range("D22").Value = 1.02435508602284
$ws.Range("E22").Value = 1.041130645914384
$ws.Range("F22").Value = 1.044876431675764
$ws.Range("I22").Value = 1.027975199674444
$ws.Range("J22").Value = 1.025565278176256
$ws.Range("K22").Value = 1.028109778548231
$ws.Range("L22").Value = 1.044820044360587
$ws.Range("M22").Value = 1.048551568217862
$ws.Range("N22").Value = 1.012476684836513
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.019064710396808
$ws.Range("D23").Value = 1.024695447258056
$ws.Range("E23").Value = 1.041662982219061
$ws.Range("F23").Value = 1.045449510364456
$ws.Range("I23").Value = 1.028056210820082
$ws.Range("J23").Value = 1.025891905106481
$ws.Range("K23").Value = 1.028379634362459
$ws.Range("L23").Value = 1.045282485791301
$ws.Range("M23").Value = 1.049054896187274
$ws.Range("N23").Value = 1.012587354907172
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.020881626145494
$ws.Range("D24").Value = 1.02603307847663
$ws.Range("E24").Value = 1.043758341118594
$ws.Range("F24").Value = 1.047705318860026
$ws.Range("I24").Value = 1.028371160915722
$ws.Range("J24").Value = 1.027174614352792
$ws.Range("K24").Value = 1.02943829364319
$ws.Range("L24").Value = 1.04710154997771
$ws.Range("M24").Value = 1.051035034576139
$ws.Range("N24").Value = 1.013021667313366
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.022985540025823
$ws.Range("D25").Value = 1.027580547424548
$ws.Range("E25").Value = 1.046189265064985
$ws.Range("F25").Value = 1.050322556696129
$ws.Range("I25").Value = 1.028728414619744
$ws.Range("J25").Value = 1.028656612484806
$ws.Range("K25").Value = 1.030659134158699
$ws.Range("L25").Value = 1.04920949784754
$ws.Range("M25").Value = 1.053330134402981
$ws.Range("N25").Value = 1.013522821509651
